$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 21 with the new incoming mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A21").Value = "Sollicitatie salesfunctie"
$ws.Range("B21").Value = "mailmind.test@zohomail.eu"
$ws.Range("C21").Value = "Hierbij mijn sollicitatie voor de salesfunctie. CV in bijlage."
$ws.Range("D21").Value = "Sollicitatie / Vacature"
$ws.Range("E21").Value = "Beste,`r`nBedankt voor je interesse in de salesfunctie bij ons bedrijf. Wij zullen je sollicitatie zorgvuldig bekijken en nemen indien nodig contact met je op. Mocht je in de tussentijd vragen hebben, aarzel dan niet om contact met ons op te nemen.`r`nMet vriendelijke groet,`r`n[Naam] `r`n[Bedrijfsnaam]"
$ws.Range("F21").Value = "2025-06-23 18:45:56"
$ws.Range("G21").Value = "Ja"

# Drop the auto row-height Excel assigns for the multi-line cell, matching the
# rest of the sheet (no explicit custom row heights anywhere else)
$ws.Rows.Item(21).AutoFit()

# Extend the conditional-formatting ranges on column D and G to cover the new row 21
$rngD = $ws.Range("D2:D21")
$rngD.FormatConditions.Item(1).ModifyAppliesToRange($rngD)

$rngG = $ws.Range("G2:G21")
$rngG.FormatConditions.Item(1).ModifyAppliesToRange($rngG)

# --- Sheet "Dashboard": re-sort category counts now that Sollicitatie / Vacature = 2 ---
$ws2 = $wb.Worksheets.Item("Dashboard")

$ws2.Range("A6").Value = "Sollicitatie / Vacature"
$ws2.Range("B6").Value = 2

$ws2.Range("A7").Value = "Offerte / Prijsaanvraag"
$ws2.Range("B7").Value = 2

$ws2.Range("A8").Value = "Productinformatie"
$ws2.Range("B8").Value = 2
